# Apply cryptocurrency price/volume updates (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.104.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5158"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3758"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07149"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8879"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("E11").Value = "  -2.59%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.868.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07546"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.307"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008479"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.97%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.127.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.005"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.092.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("E23").Value = "  -3.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.443"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.838"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.084"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.650"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.663"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09164"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05101"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.075"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.156"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7233"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.05%  "
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.091"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.488"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.074"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5276"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.479"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.288"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1465"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4614"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.966"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.563"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.16%  "
